# Updates cryptos list values (Price and Volume(1h)) to reflect the latest
# scrape, and also fixes row ordering for two coin pairs whose rank swapped
# (Filecoin/ARBITRUM at rows 34-35, and Decentraland/EnergySwap at rows 46-47).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "'27.739.66"
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").Value = "'1.904.37"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.50%  "
$ws.Range("D5").Value = "'312.77"
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("D7").Value = "'0.4982"
$ws.Range("E7").Value = "  +3.14%  "
$ws.Range("D8").Value = "'0.3788"
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").Value = "'0.07246"
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("D10").Value = "'21.14"
$ws.Range("E10").Value = "  +1.80%  "
$ws.Range("D11").Value = "'0.8988"
$ws.Range("E11").Value = "  -3.63%  "
$ws.Range("D12").Value = "'0.07628"
$ws.Range("E12").Value = "  -1.39%  "
$ws.Range("D13").Value = "'1.893.39"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").Value = "'5.455"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("D15").Value = "'91.76"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").Value = "'0.000008702"
$ws.Range("E17").Value = "  -1.82%  "
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").Value = "'27.785.64"
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("D20").Value = "'14.53"
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("D21").Value = "'5.156"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "'2.118.07"
$ws.Range("E22").Value = "  -2.80%  "
$ws.Range("E23").Value = "  -0.82%  "
$ws.Range("D24").Value = "'6.569"
$ws.Range("E24").Value = "  -0.85%  "
$ws.Range("D25").Value = "'152.93"
$ws.Range("E25").Value = "  -2.04%  "
$ws.Range("D26").Value = "'1.847"
$ws.Range("D27").Value = "'2.211"
$ws.Range("E27").Value = "  +4.23%  "
$ws.Range("D28").Value = "'18.31"
$ws.Range("E28").Value = "  -0.88%  "
$ws.Range("D29").Value = "'114.83"
$ws.Range("E29").Value = "  -2.10%  "
$ws.Range("D30").Value = "'4.869"
$ws.Range("E30").Value = "  -1.96%  "
$ws.Range("D31").Value = "'0.08922"
$ws.Range("E31").Value = "  -0.29%  "
$ws.Range("D32").Value = "'3.178"
$ws.Range("E32").Value = "  -2.11%  "
$ws.Range("D33").Value = "'0.7846"
$ws.Range("E33").Value = "  +2.41%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.229"
$ws.Range("E34").Value = "  -2.21%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "'4.786"
$ws.Range("E35").Value = "  +2.58%  "
$ws.Range("D36").Value = "'2.614"
$ws.Range("E36").Value = "  +2.61%  "
$ws.Range("D37").Value = "'0.02073"
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("D38").Value = "'3.058"
$ws.Range("E38").Value = "  +2.01%  "
$ws.Range("D39").Value = "'1.090"
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("D40").Value = "'0.5506"
$ws.Range("E40").Value = "  +0.44%  "
$ws.Range("D41").Value = "'0.05291"
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("D42").Value = "'6.748"
$ws.Range("E42").Value = "  -2.89%  "
$ws.Range("D43").Value = "'114.24"
$ws.Range("E43").Value = "  +3.77%  "
$ws.Range("D44").Value = "'8.459"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("E45").Value = "  -1.02%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'10.48"
$ws.Range("E46").Value = "  -1.59%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.4778"
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("D48").Value = "'1.001"
$ws.Range("E48").Value = "  -0.46%  "
$ws.Range("D49").Value = "'1.630"
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("D50").Value = "'67.00"
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("D51").Value = "'0.06008"
$ws.Range("E51").Value = "  -1.17%  "
